# Update "想去人数" (column F) and "最低票价" (column G) figures on the
# "展览" (sheet1), "演出" (sheet2) and "全部类型" (sheet4) worksheets to match
# the regenerated gh-pages data snapshot at commit 456a3b4.

$sheet1Changes = @(
    @{ Row = 2; F = 200; G = $null },
    @{ Row = 3; F = 1373; G = $null },
    @{ Row = 4; F = 19300; G = 60 },
    @{ Row = 6; F = 295; G = $null },
    @{ Row = 8; F = 4; G = $null },
    @{ Row = 9; F = 7310; G = $null },
    @{ Row = 10; F = 471; G = $null },
    @{ Row = 11; F = 713; G = $null },
    @{ Row = 12; F = 236; G = $null },
    @{ Row = 13; F = 27; G = $null },
    @{ Row = 17; F = 178; G = $null },
    @{ Row = 18; F = 1323; G = $null },
    @{ Row = 19; F = 342; G = $null },
    @{ Row = 20; F = 65; G = $null },
    @{ Row = 21; F = 671; G = $null },
    @{ Row = 23; F = 45; G = $null },
    @{ Row = 24; F = 54; G = $null },
    @{ Row = 25; F = 298; G = $null },
    @{ Row = 26; F = 1057; G = $null },
    @{ Row = 27; F = 19; G = $null },
    @{ Row = 28; F = 4; G = $null },
    @{ Row = 29; F = 154; G = $null },
    @{ Row = 30; F = 5219; G = $null },
    @{ Row = 31; F = 549; G = $null },
    @{ Row = 32; F = 41; G = $null },
    @{ Row = 33; F = 136; G = $null },
    @{ Row = 35; F = 82; G = $null },
    @{ Row = 36; F = 12409; G = $null },
    @{ Row = 37; F = 1312; G = $null },
    @{ Row = 38; F = 46; G = $null },
    @{ Row = 42; F = 325; G = $null },
    @{ Row = 43; F = 3963; G = $null },
    @{ Row = 44; F = 315; G = $null },
    @{ Row = 45; F = 95; G = $null }
)

$sheet2Changes = @(
    @{ Row = 3; F = 32; G = $null }
)

$sheet4Changes = @(
    @{ Row = 2; F = 200; G = $null },
    @{ Row = 3; F = 1373; G = $null },
    @{ Row = 4; F = 19300; G = 60 },
    @{ Row = 5; F = 767; G = $null },
    @{ Row = 6; F = 295; G = $null },
    @{ Row = 7; F = 1084; G = $null },
    @{ Row = 8; F = 4; G = $null },
    @{ Row = 9; F = 7310; G = $null },
    @{ Row = 10; F = 471; G = $null },
    @{ Row = 11; F = 713; G = $null },
    @{ Row = 12; F = 236; G = $null },
    @{ Row = 13; F = 27; G = $null },
    @{ Row = 14; F = 141; G = $null },
    @{ Row = 15; F = 90; G = $null },
    @{ Row = 18; F = 1323; G = $null },
    @{ Row = 19; F = 342; G = $null },
    @{ Row = 21; F = 671; G = $null },
    @{ Row = 22; F = 42; G = $null },
    @{ Row = 23; F = 45; G = $null },
    @{ Row = 24; F = 54; G = $null },
    @{ Row = 25; F = 298; G = $null },
    @{ Row = 26; F = 1057; G = $null },
    @{ Row = 27; F = 19; G = $null },
    @{ Row = 28; F = 4; G = $null },
    @{ Row = 29; F = 154; G = $null },
    @{ Row = 30; F = 5219; G = $null },
    @{ Row = 31; F = 549; G = $null },
    @{ Row = 33; F = 41; G = $null },
    @{ Row = 34; F = 32; G = $null },
    @{ Row = 35; F = 136; G = $null },
    @{ Row = 38; F = 12409; G = $null },
    @{ Row = 39; F = 1312; G = $null },
    @{ Row = 40; F = 46; G = $null },
    @{ Row = 43; F = 246; G = $null },
    @{ Row = 45; F = 3963; G = $null },
    @{ Row = 47; F = 95; G = $null }
)


$wb = $excel.ActiveWorkbook

function Apply-Changes($ws, $changes) {
    foreach ($chg in $changes) {
        if ($null -ne $chg.F) {
            $ws.Cells.Item($chg.Row, 6).Value = $chg.F
        }
        if ($null -ne $chg.G) {
            $ws.Cells.Item($chg.Row, 7).Value = $chg.G
        }
    }
}

$wsExhibition = $wb.Worksheets.Item("展览")
Apply-Changes $wsExhibition $sheet1Changes

$wsPerformance = $wb.Worksheets.Item("演出")
Apply-Changes $wsPerformance $sheet2Changes

$wsAllTypes = $wb.Worksheets.Item("全部类型")
Apply-Changes $wsAllTypes $sheet4Changes
